# Import ADJ column from excel budget file.
# Add three new header columns (ADJ1, ADJ2, ADJ3) after the existing
# "Project" header column (R1), at S1:U1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

$ws.Range("S1").Value = "ADJ1"
$ws.Range("T1").Value = "ADJ2"
$ws.Range("U1").Value = "ADJ3"

$ws.Range("S1:U1").Select()
